# Re-process the metadata with the newly curated dimensions.
# Column C = municipio-nombre, Column E = n-de-habitaciones-de-la-vivienda

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Column C (municipio-nombre) --------------------------------------
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"

# -- Column E (n-de-habitaciones-de-la-vivienda) -----------------------
$ws.Range("E2").Value = "iaest-measure:n-de-habitaciones-de-la-vivienda"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("E5").Clear()
